$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-26 Thursday" "2024-09-27 Friday"

Replace-Text "356÷2=178, 0" "643÷2=321, 1"
Replace-Text "783÷2=391, 1" "347÷4=86, 3"
Replace-Text "244÷2=122, 0" "222÷5=44, 2"
Replace-Text "178÷6=29, 4" "789÷8=98, 5"
Replace-Text "170÷2=85, 0" "539÷3=179, 2"

Replace-Text "190÷2=95, 0" "482÷2=241, 0"
Replace-Text "945÷2=472, 1" "305÷5=61, 0"
Replace-Text "195÷4=48, 3" "354÷3=118, 0"
Replace-Text "104÷3=34, 2" "235÷2=117, 1"
Replace-Text "334÷9=37, 1" "300÷2=150, 0"

Replace-Text "783÷4=195, 3" "286÷5=57, 1"
Replace-Text "579÷5=115, 4" "978÷4=244, 2"
Replace-Text "718÷9=79, 7" "269÷8=33, 5"
Replace-Text "221÷3=73, 2" "707÷2=353, 1"
Replace-Text "335÷9=37, 2" "329÷8=41, 1"

Replace-Text "139÷9=15, 4" "163÷9=18, 1"
Replace-Text "520÷7=74, 2" "783÷6=130, 3"
Replace-Text "874÷4=218, 2" "966÷9=107, 3"
Replace-Text "699÷3=233, 0" "689÷3=229, 2"
Replace-Text "918÷8=114, 6" "860÷5=172, 0"

Replace-Text "689÷5=137, 4" "650÷6=108, 2"
Replace-Text "243÷5=48, 3" "747÷8=93, 3"
Replace-Text "628÷5=125, 3" "928÷8=116, 0"
Replace-Text "757÷7=108, 1" "490÷5=98, 0"
Replace-Text "337÷7=48, 1" "471÷3=157, 0"
